# Apply updated cryptocurrency price/volume data to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in this sheet are stored as text (inline strings), including
# price figures that look numeric (e.g. "218.60", "1.861.38"). Force the cells
# to Text format before assigning so Excel does not auto-convert them to numbers
# or drop significant trailing/leading characters (e.g. trailing zeros).
$cellValues = [ordered]@{
    "D2" = "28.320.32"
    "E2" = "  +4.17%  "
    "D3" = "1.716.68"
    "E3" = "  +1.76%  "
    "E4" = "  -0.10%  "
    "D5" = "218.60"
    "E5" = "  +1.40%  "
    "E6" = "  +0.25%  "
    "E7" = "  -0.07%  "
    "D8" = "23.90"
    "E8" = "  +3.22%  "
    "E9" = "  +2.06%  "
    "D10" = "0.0633"
    "E10" = "  +0.95%  "
    "D11" = "0.0891"
    "E11" = "  +0.09%  "
    "D12" = "1.958.47"
    "E12" = "  +1.54%  "
    "D13" = "1.710.18"
    "E13" = "  +1.15%  "
    "D14" = "4.22"
    "E14" = "  +0.37%  "
    "D15" = "0.561"
    "E15" = "  +1.03%  "
    "D16" = "67.38"
    "E16" = "  +0.36%  "
    "D17" = "28.321.42"
    "E17" = "  +4.02%  "
    "D18" = "248.86"
    "E18" = "  +5.20%  "
    "D19" = "0.0₃0747"
    "E19" = "  +0.49%  "
    "D20" = "7.78"
    "E20" = "  -3.58%  "
    "E21" = "  +0.05%  "
    "E22" = "  +0.42%  "
    "D23" = "9.60"
    "E23" = "  -0.06%  "
    "E24" = "  -1.56%  "
    "D25" = "147.52"
    "E25" = "  +0.23%  "
    "E26" = "  +0.90%  "
    "D27" = "16.52"
    "E27" = "  +0.32%  "
    "E28" = "  +0.58%  "
    "E29" = "  -0.16%  "
    "E30" = "  +1.43%  "
    "E31" = "  +2.47%  "
    "E32" = "  +0.37%  "
    "D33" = "1.476.32"
    "E33" = "  -4.30%  "
    "E34" = "  -1.13%  "
    "D35" = "1.63"
    "E35" = "  -2.15%  "
    "D36" = "0.973"
    "E36" = "  +2.97%  "
    "D37" = "2.40"
    "E37" = "  +0.42%  "
    "E38" = "  -1.35%  "
    "E39" = "  +0.77%  "
    "E40" = "  -0.63%  "
    "D41" = "69.51"
    "E41" = "  +0.32%  "
    "E42" = "  -0.04%  "
    "D43" = "5.64"
    "E43" = "  -1.96%  "
    "B44" = "MXToken"
    "C44" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D44" = "2.28"
    "E44" = "  +0.54%  "
    "B45" = "RocketPoolETH"
    "C45" = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
    "D45" = "1.862.59"
    "E45" = "  +1.27%  "
    "D46" = "0.804"
    "E46" = "  +1.28%  "
    "E47" = "  +6.45%  "
    "D48" = "89.92"
    "E48" = "  -0.69%  "
    "D49" = "0.0₆0112"
    "E49" = "  +0.06%  "
    "D50" = "8.08"
    "E50" = "  -3.22%  "
    "E51" = "  -1.19%  "
}

foreach ($cellRef in $cellValues.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $cellValues[$cellRef]
}
